# Applies the "feat: add 2022-Q4 data" edit:
#  1. Insert a new worksheet "2022-Q4" right after "总计" and before "2022-Q3",
#     populated with the new quarter's per-fund holdings table.
#  2. Insert a new summary row for "2022-Q4" at the top of the data in "总计"
#     (row 2), pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q4" worksheet, placed right after "总计".
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item(2)          # currently "2022-Q3"
$q4 = $wb.Worksheets.Add($afterSheet)         # inserted before it -> right after 总计
$q4.Name = "2022-Q4"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = [char]([int][char]'B' + $i)
    $cell = $q4.Range($col + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @(0,"010387","易方达医药生物股票A","20.77","91.88","3.65","0.7581",10),
    @(1,"007718","中银创新医疗混合A","11.61","80.72","3.37","0.3913",10),
    @(2,"012706","中银核心精选混合A","3.65","93.68","9.65","0.3522",1),
    @(3,"009877","中银内核驱动股票A","2.58","92.75","9.59","0.2474",2),
    @(4,"010388","易方达医药生物股票C","6.05","91.88","3.65","0.2208",10),
    @(5,"862001","光大阳光香港精选混合（QDII）A 人民币","3.90","92.65","4.53","0.1767",8),
    @(6,"862011","光大阳光香港精选混合（QDII）A 美元","3.90","92.65","4.53","0.1767",8),
    @(7,"862012","光大阳光香港精选混合（QDII）C 人民币","3.90","92.65","4.53","0.1767",8),
    @(8,"009898","民生加银医药健康股票A","2.62","89.21","5.12","0.1341",2),
    @(9,"010500","中银创新医疗混合C","3.82","80.72","3.37","0.1287",10),
    @(10,"580008","东吴新产业精选股票A","4.38","91.75","2.78","0.1218",7),
    @(11,"011470","东吴新产业精选混合C","4.38","91.75","2.78","0.1218",7),
    @(12,"007182","万家沪港深蓝筹混合A","3.44","92.13","3.49","0.1201",9),
    @(13,"011157","弘毅远方港股通智选领航混合A","2.84","90.73","4.03","0.1145",7),
    @(14,"005805","华泰柏瑞医疗健康混合A","4.83","80.25","2.17","0.1048",10),
    @(15,"009353","浙商科技创新一个月滚动持有混合A","1.49","90.77","6.94","0.1034",4),
    @(16,"012584","南方中国新兴经济9个月持有期混合（QDII）A","3.08","83.21","2.89","0.0890",9),
    @(17,"013009","万家港股通精选混合A","2.56","84.90","3.30","0.0845",6),
    @(18,"513120","广发中证香港创新药（QDII-ETF）","1.85","98.69","4.16","0.0770",8),
    @(19,"009354","浙商科技创新一个月滚动持有混合C","1.02","90.77","6.94","0.0708",4),
    @(20,"006603","嘉实互融精选股票","0.63","91.84","9.39","0.0592",2),
    @(21,"015373","浙商智选新兴产业混合A","0.72","92.09","6.66","0.0480",2),
    @(22,"005029","中银产业精选混合A","0.51","93.13","7.50","0.0382",6),
    @(23,"011158","弘毅远方港股通智选领航混合C","0.84","90.73","4.03","0.0339",7),
    @(24,"007183","万家沪港深蓝筹混合C","0.86","92.13","3.49","0.0300",9),
    @(25,"013010","万家港股通精选混合C","0.77","84.90","3.30","0.0254",6),
    @(26,"006072","民生加银创新成长混合A","0.40","91.73","5.67","0.0227",1),
    @(27,"005520","国投瑞银创新医疗混合","0.49","92.50","4.13","0.0202",5),
    @(28,"015374","浙商智选新兴产业混合C","0.29","92.09","6.66","0.0193",2),
    @(29,"011453","华泰柏瑞医疗健康混合C","0.67","80.25","2.17","0.0145",10),
    @(30,"860008","光大阳光生活18个月持有期混合A","0.33","90.14","3.88","0.0128",2),
    @(31,"012707","中银核心精选混合C","0.13","93.68","9.65","0.0125",1),
    @(32,"012315","创金合信港股通成长股票A","0.12","89.18","9.89","0.0119",1),
    @(33,"012316","创金合信港股通成长股票C","0.11","89.18","9.89","0.0109",1),
    @(34,"008861","西部利得港股通新机遇灵活配置混合A","0.25","87.69","3.65","0.0091",5),
    @(35,"860060","光大阳光生活18个月持有期混合B","0.18","90.14","3.88","0.0070",2),
    @(36,"013182","安信港股通精选混合C","0.12","69.28","4.54","0.0054",3),
    @(37,"010093","西部利得港股通新机遇灵活配置混合C","0.12","87.69","3.65","0.0044",5),
    @(38,"012585","南方中国新兴经济9个月持有期混合（QDII）C","0.12","83.21","2.89","0.0035",9),
    @(39,"860061","光大阳光生活18个月持有期混合C","0.08","90.14","3.88","0.0031",2),
    @(40,"012600","中银内核驱动股票C","0.02","92.75","9.59","0.0019",2),
    @(41,"013181","安信港股通精选混合A","0.02","69.28","4.54","0.0009",3),
    @(42,"014929","民生加银创新成长混合C","0.01","91.73","5.67","0.0006",1),
    @(43,"014758","民生加银医药健康股票C","0.01","89.21","5.12","0.0005",2),
    @(44,"005030","中银产业精选混合C","0.00","93.13","7.50",0,6)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $aCell = $q4.Range("A" + $r)
    $aCell.Value = $row[0]
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $q4.Range("B" + $r).Value = $row[1]
    $q4.Range("C" + $r).Value = $row[2]
    $q4.Range("D" + $r).Value = $row[3]
    $q4.Range("E" + $r).Value = $row[4]
    $q4.Range("F" + $r).Value = $row[5]
    $q4.Range("G" + $r).Value = $row[6]
    $q4.Range("H" + $r).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: shift rows 2-7 down to 3-8, then write
#    the new 2022-Q4 summary row at row 2.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

for ($r = 7; $r -ge 2; $r--) {
    $newR = $r + 1
    $total.Range("A" + $r).Copy($total.Range("A" + $newR))
    $total.Range("B" + $r).Copy($total.Range("B" + $newR))
    $total.Range("C" + $r).Copy($total.Range("C" + $newR))
    $total.Range("D" + $r).Copy($total.Range("D" + $newR))
}

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 45
$total.Range("D2").Value = 4.17
